$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.231.59'
$ws.Range("E2").Value = '  +2.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.086.16'
$ws.Range("E3").Value = '  +3.07%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.73'
$ws.Range("E5").Value = '  +3.00%  '

$ws.Range("E6").Value = '  +1.53%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.71'
$ws.Range("E8").Value = '  +25.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.63'
$ws.Range("E9").Value = '  +3.15%  '

$ws.Range("E10").Value = '  +6.57%  '

$ws.Range("E11").Value = '  +5.41%  '

$ws.Range("E12").Value = '  +8.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.02'
$ws.Range("E13").Value = '  +6.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.391.49'
$ws.Range("E14").Value = '  +3.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.829'
$ws.Range("E15").Value = '  +4.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.092.25'
$ws.Range("E16").Value = '  +3.55%  '

$ws.Range("E17").Value = '  +7.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.159.01'
$ws.Range("E18").Value = '  +2.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.72'
$ws.Range("E19").Value = '  +3.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.52'
$ws.Range("E20").Value = '  +17.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0846'
$ws.Range("E21").Value = '  +5.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.16'
$ws.Range("E22").Value = '  +2.89%  '

$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("E25").Value = '  +1.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.19'
$ws.Range("E26").Value = '  +2.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.18'
$ws.Range("E27").Value = '  +6.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.66'
$ws.Range("E28").Value = '  +5.10%  '

$ws.Range("E29").Value = '  +5.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.123'
$ws.Range("E30").Value = '  +2.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("E31").Value = '  +32.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.27'
$ws.Range("E32").Value = '  +4.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.52'
$ws.Range("E33").Value = '  +5.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0619'
$ws.Range("E34").Value = '  +8.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0898'
$ws.Range("E35").Value = '  +1.77%  '

$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("E37").Value = '  +6.31%  '

$ws.Range("E38").Value = '  -3.32%  '

$ws.Range("E39").Value = '  +5.15%  '

$ws.Range("B40").Value = 'FTXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.17'
$ws.Range("E40").Value = '  +174.82%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.35'
$ws.Range("E41").Value = '  +4.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '18.18'
$ws.Range("E42").Value = '  +18.17%  '

$ws.Range("E44").Value = '  +7.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '98.41'
$ws.Range("E45").Value = '  +3.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0950'
$ws.Range("E46").Value = '  +16.97%  '

$ws.Range("E47").Value = '  +0.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.323.76'
$ws.Range("E48").Value = '  +2.03%  '

$ws.Range("E49").Value = '  +5.54%  '

$ws.Range("E50").Value = '  +8.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.99'
$ws.Range("E51").Value = '  +16.02%  '
